$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 85 (swap with former row 86 content, id stays 83)
$ws.Range("A85").Value = 83
$ws.Range("B85").Value = 6992620
$ws.Range("C85").Value = "Thailand Premier League"
$ws.Range("D85").Value = "Thailand Premier League"
$ws.Range("E85").Value = 45261.375
$ws.Range("F85").Value = "Uthai Thani FC"
$ws.Range("G85").Value = "Sukhothai FC"
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = "D"
$ws.Range("K85").Value = 1.95
$ws.Range("L85").Value = 3.5
$ws.Range("M85").Value = 3.4
$ws.Range("N85").Value = 2.1
$ws.Range("O85").Value = 3.4
$ws.Range("P85").Value = 3
$ws.Range("Q85").Value = -0.25
$ws.Range("R85").Value = 1.875
$ws.Range("S85").Value = 1.925
$ws.Range("T85").Value = 2.75
$ws.Range("U85").Value = 1.8
$ws.Range("V85").Value = 2
$ws.Range("W85").Value = -1
$ws.Range("X85").Value = 2.4
$ws.Range("Y85").Value = -1
$ws.Range("Z85").Value = -0.5
$ws.Range("AA85").Value = 0.4625
$ws.Range("AB85").Value = -1
$ws.Range("AC85").Value = 1

# Row 86 (swap with former row 85 content, id stays 84)
$ws.Range("A86").Value = 84
$ws.Range("B86").Value = 6992623
$ws.Range("C86").Value = "Thailand Premier League"
$ws.Range("D86").Value = "Thailand Premier League"
$ws.Range("E86").Value = 45261.375
$ws.Range("F86").Value = "Ratchaburi FC"
$ws.Range("G86").Value = "Chiangrai Utd"
$ws.Range("H86").Value = 3
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = "H"
$ws.Range("K86").Value = 1.7
$ws.Range("L86").Value = 3.75
$ws.Range("M86").Value = 4.2
$ws.Range("N86").Value = 1.7
$ws.Range("O86").Value = 3.75
$ws.Range("P86").Value = 4.333
$ws.Range("Q86").Value = -0.75
$ws.Range("R86").Value = 1.925
$ws.Range("S86").Value = 1.875
$ws.Range("T86").Value = 2.5
$ws.Range("U86").Value = 1.85
$ws.Range("V86").Value = 1.95
$ws.Range("W86").Value = 0.7
$ws.Range("X86").Value = -1
$ws.Range("Y86").Value = -1
$ws.Range("Z86").Value = 0.925
$ws.Range("AA86").Value = -1
$ws.Range("AB86").Value = 0.8500000000000001
$ws.Range("AC86").Value = -1

# Row 117 (swap with former row 118 content, id stays 115)
$ws.Range("A117").Value = 115
$ws.Range("B117").Value = 7329293
$ws.Range("C117").Value = "Thailand Premier League"
$ws.Range("D117").Value = "Thailand Premier League"
$ws.Range("E117").Value = 45288.375
$ws.Range("F117").Value = "Chonburi"
$ws.Range("G117").Value = "Bangkok United"
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = "D"
$ws.Range("K117").Value = 3.6
$ws.Range("L117").Value = 3.5
$ws.Range("M117").Value = 1.85
$ws.Range("N117").Value = 4.5
$ws.Range("O117").Value = 4
$ws.Range("P117").Value = 1.615
$ws.Range("Q117").Value = 0.75
$ws.Range("R117").Value = 1.975
$ws.Range("S117").Value = 1.825
$ws.Range("T117").Value = 3
$ws.Range("U117").Value = 1.85
$ws.Range("V117").Value = 1.95
$ws.Range("W117").Value = -1
$ws.Range("X117").Value = 3
$ws.Range("Y117").Value = -1
$ws.Range("Z117").Value = 0.9750000000000001
$ws.Range("AA117").Value = -1
$ws.Range("AB117").Value = -1
$ws.Range("AC117").Value = 0.95

# Row 118 (swap with former row 117 content, id stays 116)
$ws.Range("A118").Value = 116
$ws.Range("B118").Value = 7485127
$ws.Range("C118").Value = "Thailand Premier League"
$ws.Range("D118").Value = "Thailand Premier League"
$ws.Range("E118").Value = 45288.375
$ws.Range("F118").Value = "BG Pathum United"
$ws.Range("G118").Value = "Chiangrai Utd"
$ws.Range("H118").Value = 2
$ws.Range("I118").Value = 2
$ws.Range("J118").Value = "D"
$ws.Range("K118").Value = 1.5
$ws.Range("L118").Value = 4
$ws.Range("M118").Value = 5.75
$ws.Range("N118").Value = 1.363
$ws.Range("O118").Value = 4.5
$ws.Range("P118").Value = 6.5
$ws.Range("Q118").Value = -1.25
$ws.Range("R118").Value = 1.85
$ws.Range("S118").Value = 1.95
$ws.Range("T118").Value = 3
$ws.Range("U118").Value = 1.825
$ws.Range("V118").Value = 1.975
$ws.Range("W118").Value = -1
$ws.Range("X118").Value = 3.5
$ws.Range("Y118").Value = -1
$ws.Range("Z118").Value = -1
$ws.Range("AA118").Value = 0.95
$ws.Range("AB118").Value = 0.825
$ws.Range("AC118").Value = -1

# Row 155 updates (new H/I/J + odds changes)
$ws.Range("H155").Value = 2
$ws.Range("I155").Value = 1
$ws.Range("J155").Value = "H"
$ws.Range("K155").Value = 1.95
$ws.Range("L155").Value = 3.6
$ws.Range("M155").Value = 3.2
$ws.Range("N155").Value = 1.727
$ws.Range("O155").Value = 3.8
$ws.Range("P155").Value = 3.75
$ws.Range("Q155").Value = -0.75
$ws.Range("R155").Value = 1.975
$ws.Range("S155").Value = 1.825
$ws.Range("T155").Value = 3
$ws.Range("U155").Value = 1.925
$ws.Range("V155").Value = 1.875
$ws.Range("W155").Value = 0.7270000000000001
$ws.Range("X155").Value = -1
$ws.Range("Y155").Value = -1
$ws.Range("Z155").Value = 0.4875
$ws.Range("AA155").Value = -0.5
$ws.Range("AB155").Value = 0
$ws.Range("AC155").Value = 0

# Row 156 updates (new H/I/J + odds changes)
$ws.Range("H156").Value = 3
$ws.Range("I156").Value = 2
$ws.Range("J156").Value = "H"
$ws.Range("K156").Value = 1.833
$ws.Range("L156").Value = 3.6
$ws.Range("M156").Value = 3.5
$ws.Range("N156").Value = 1.75
$ws.Range("O156").Value = 3.6
$ws.Range("P156").Value = 3.8
$ws.Range("Q156").Value = -0.75
$ws.Range("R156").Value = 2
$ws.Range("S156").Value = 1.8
$ws.Range("T156").Value = 2.5
$ws.Range("U156").Value = 1.8
$ws.Range("V156").Value = 2
$ws.Range("W156").Value = 0.75
$ws.Range("X156").Value = -1
$ws.Range("Y156").Value = -1
$ws.Range("Z156").Value = 0.5
$ws.Range("AA156").Value = -0.5
$ws.Range("AB156").Value = 0.8
$ws.Range("AC156").Value = -1

# Row 157 updates (new H/I/J + odds changes)
$ws.Range("H157").Value = 1
$ws.Range("I157").Value = 1
$ws.Range("J157").Value = "D"
$ws.Range("K157").Value = 1.95
$ws.Range("L157").Value = 3.5
$ws.Range("M157").Value = 3.25
$ws.Range("N157").Value = 2.15
$ws.Range("O157").Value = 3.4
$ws.Range("P157").Value = 2.8
$ws.Range("Q157").Value = -0.25
$ws.Range("R157").Value = 2
$ws.Range("S157").Value = 1.8
$ws.Range("T157").Value = 2.75
$ws.Range("U157").Value = 1.95
$ws.Range("V157").Value = 1.85
$ws.Range("W157").Value = -1
$ws.Range("X157").Value = 2.4
$ws.Range("Y157").Value = -1
$ws.Range("Z157").Value = -0.5
$ws.Range("AA157").Value = 0.4
$ws.Range("AB157").Value = -1
$ws.Range("AC157").Value = 0.8500000000000001

# Row 158 updates (new H/I/J + odds changes)
$ws.Range("H158").Value = 2
$ws.Range("I158").Value = 2
$ws.Range("J158").Value = "D"
$ws.Range("K158").Value = 2.875
$ws.Range("L158").Value = 3.4
$ws.Range("M158").Value = 2.15
$ws.Range("N158").Value = 3.3
$ws.Range("O158").Value = 3.3
$ws.Range("P158").Value = 2
$ws.Range("Q158").Value = 0.5
$ws.Range("R158").Value = 1.725
$ws.Range("S158").Value = 1.975
$ws.Range("T158").Value = 2.5
$ws.Range("U158").Value = 1.875
$ws.Range("V158").Value = 1.925
$ws.Range("W158").Value = -1
$ws.Range("X158").Value = 2.3
$ws.Range("Y158").Value = -1
$ws.Range("Z158").Value = 0.7250000000000001
$ws.Range("AA158").Value = -1
$ws.Range("AB158").Value = 0.875
$ws.Range("AC158").Value = -1
# New rows 159-166

# Row 159
$ws.Range("A158").Copy($ws.Range("A159"))
$ws.Range("E158").Copy($ws.Range("E159"))
$ws.Range("A159").Value = 157
$ws.Range("B159").Value = 6992334
$ws.Range("C159").Value = "Thailand Premier League"
$ws.Range("D159").Value = "Thailand Premier League"
$ws.Range("E159").Value = 45359.375
$ws.Range("F159").Value = "Prachuap FC"
$ws.Range("G159").Value = "Police Tero FC"
$ws.Range("K159").Value = 1.533
$ws.Range("L159").Value = 4.2
$ws.Range("M159").Value = 4.5
$ws.Range("N159").Value = 1.5
$ws.Range("O159").Value = 4.2
$ws.Range("P159").Value = 4.75
$ws.Range("Q159").Value = -1
$ws.Range("R159").Value = 1.825
$ws.Range("S159").Value = 1.975
$ws.Range("T159").Value = 3
$ws.Range("U159").Value = 2
$ws.Range("V159").Value = 1.8
$ws.Range("W159").Value = 0
$ws.Range("X159").Value = 0
$ws.Range("Y159").Value = 0
$ws.Range("Z159").Value = 0
$ws.Range("AA159").Value = 0

# Row 160
$ws.Range("A158").Copy($ws.Range("A160"))
$ws.Range("E158").Copy($ws.Range("E160"))
$ws.Range("A160").Value = 158
$ws.Range("B160").Value = 6992683
$ws.Range("C160").Value = "Thailand Premier League"
$ws.Range("D160").Value = "Thailand Premier League"
$ws.Range("E160").Value = 45360.33333333334
$ws.Range("F160").Value = "Chiangrai Utd"
$ws.Range("G160").Value = "Bangkok United"
$ws.Range("K160").Value = 5.5
$ws.Range("L160").Value = 3.75
$ws.Range("M160").Value = 1.5
$ws.Range("N160").Value = 5.5
$ws.Range("O160").Value = 3.75
$ws.Range("P160").Value = 1.5
$ws.Range("Q160").Value = 1
$ws.Range("R160").Value = 1.95
$ws.Range("S160").Value = 1.85
$ws.Range("T160").Value = 2.5
$ws.Range("U160").Value = 1.925
$ws.Range("V160").Value = 1.875
$ws.Range("W160").Value = 0
$ws.Range("X160").Value = 0
$ws.Range("Y160").Value = 0
$ws.Range("Z160").Value = 0
$ws.Range("AA160").Value = 0

# Row 161
$ws.Range("A158").Copy($ws.Range("A161"))
$ws.Range("E158").Copy($ws.Range("E161"))
$ws.Range("A161").Value = 159
$ws.Range("B161").Value = 6992682
$ws.Range("C161").Value = "Thailand Premier League"
$ws.Range("D161").Value = "Thailand Premier League"
$ws.Range("E161").Value = 45360.35416666666
$ws.Range("F161").Value = "BG Pathum United"
$ws.Range("G161").Value = "Sukhothai FC"
$ws.Range("K161").Value = 1.3
$ws.Range("L161").Value = 4.75
$ws.Range("M161").Value = 7.5
$ws.Range("N161").Value = 1.363
$ws.Range("O161").Value = 4.5
$ws.Range("P161").Value = 6
$ws.Range("Q161").Value = -1.25
$ws.Range("R161").Value = 1.825
$ws.Range("S161").Value = 1.975
$ws.Range("T161").Value = 2.75
$ws.Range("U161").Value = 1.825
$ws.Range("V161").Value = 1.975
$ws.Range("W161").Value = 0
$ws.Range("X161").Value = 0
$ws.Range("Y161").Value = 0
$ws.Range("Z161").Value = 0
$ws.Range("AA161").Value = 0

# Row 162
$ws.Range("A158").Copy($ws.Range("A162"))
$ws.Range("E158").Copy($ws.Range("E162"))
$ws.Range("A162").Value = 160
$ws.Range("B162").Value = 6992686
$ws.Range("C162").Value = "Thailand Premier League"
$ws.Range("D162").Value = "Thailand Premier League"
$ws.Range("E162").Value = 45360.375
$ws.Range("F162").Value = "Ratchaburi FC"
$ws.Range("G162").Value = "Lamphun Warrior FC"
$ws.Range("K162").Value = 1.833
$ws.Range("L162").Value = 3.4
$ws.Range("M162").Value = 3.6
$ws.Range("N162").Value = 1.833
$ws.Range("O162").Value = 3.3
$ws.Range("P162").Value = 3.75
$ws.Range("Q162").Value = -0.5
$ws.Range("R162").Value = 1.875
$ws.Range("S162").Value = 1.925
$ws.Range("T162").Value = 2.5
$ws.Range("U162").Value = 2
$ws.Range("V162").Value = 1.8
$ws.Range("W162").Value = 0
$ws.Range("X162").Value = 0
$ws.Range("Y162").Value = 0
$ws.Range("Z162").Value = 0
$ws.Range("AA162").Value = 0

# Row 163
$ws.Range("A158").Copy($ws.Range("A163"))
$ws.Range("E158").Copy($ws.Range("E163"))
$ws.Range("A163").Value = 161
$ws.Range("B163").Value = 6992684
$ws.Range("C163").Value = "Thailand Premier League"
$ws.Range("D163").Value = "Thailand Premier League"
$ws.Range("E163").Value = 45360.41666666666
$ws.Range("F163").Value = "Port FC"
$ws.Range("G163").Value = "Uthai Thani FC"
$ws.Range("K163").Value = 1.4
$ws.Range("L163").Value = 4.5
$ws.Range("M163").Value = 6
$ws.Range("N163").Value = 1.4
$ws.Range("O163").Value = 4.5
$ws.Range("P163").Value = 6
$ws.Range("Q163").Value = -1.5
$ws.Range("R163").Value = 2
$ws.Range("S163").Value = 1.8
$ws.Range("T163").Value = 3.25
$ws.Range("U163").Value = 1.925
$ws.Range("V163").Value = 1.875
$ws.Range("W163").Value = 0
$ws.Range("X163").Value = 0
$ws.Range("Y163").Value = 0
$ws.Range("Z163").Value = 0
$ws.Range("AA163").Value = 0

# Row 164
$ws.Range("A158").Copy($ws.Range("A164"))
$ws.Range("E158").Copy($ws.Range("E164"))
$ws.Range("A164").Value = 162
$ws.Range("B164").Value = 6992685
$ws.Range("C164").Value = "Thailand Premier League"
$ws.Range("D164").Value = "Thailand Premier League"
$ws.Range("E164").Value = 45361.3125
$ws.Range("F164").Value = "Nakhon Pathom FC"
$ws.Range("G164").Value = "Trat FC"
$ws.Range("K164").Value = 2.05
$ws.Range("L164").Value = 3.4
$ws.Range("M164").Value = 3
$ws.Range("N164").Value = 2
$ws.Range("O164").Value = 3.4
$ws.Range("P164").Value = 3.1
$ws.Range("Q164").Value = -0.25
$ws.Range("R164").Value = 1.8
$ws.Range("S164").Value = 2
$ws.Range("T164").Value = 2.75
$ws.Range("U164").Value = 1.925
$ws.Range("V164").Value = 1.875
$ws.Range("W164").Value = 0
$ws.Range("X164").Value = 0
$ws.Range("Y164").Value = 0
$ws.Range("Z164").Value = 0
$ws.Range("AA164").Value = 0

# Row 165
$ws.Range("A158").Copy($ws.Range("A165"))
$ws.Range("E158").Copy($ws.Range("E165"))
$ws.Range("A165").Value = 163
$ws.Range("B165").Value = 6992681
$ws.Range("C165").Value = "Thailand Premier League"
$ws.Range("D165").Value = "Thailand Premier League"
$ws.Range("E165").Value = 45361.33333333334
$ws.Range("F165").Value = "Chonburi"
$ws.Range("G165").Value = "Buriram United"
$ws.Range("K165").Value = 4.5
$ws.Range("L165").Value = 4
$ws.Range("M165").Value = 1.571
$ws.Range("N165").Value = 4.5
$ws.Range("O165").Value = 3.8
$ws.Range("P165").Value = 1.6
$ws.Range("Q165").Value = 1
$ws.Range("R165").Value = 1.825
$ws.Range("S165").Value = 1.975
$ws.Range("T165").Value = 3
$ws.Range("U165").Value = 1.875
$ws.Range("V165").Value = 1.925
$ws.Range("W165").Value = 0
$ws.Range("X165").Value = 0
$ws.Range("Y165").Value = 0
$ws.Range("Z165").Value = 0
$ws.Range("AA165").Value = 0

# Row 166
$ws.Range("A158").Copy($ws.Range("A166"))
$ws.Range("E158").Copy($ws.Range("E166"))
$ws.Range("A166").Value = 164
$ws.Range("B166").Value = 6992687
$ws.Range("C166").Value = "Thailand Premier League"
$ws.Range("D166").Value = "Thailand Premier League"
$ws.Range("E166").Value = 45361.375
$ws.Range("F166").Value = "Muang Thong United"
$ws.Range("G166").Value = "Khonkaen United"
$ws.Range("K166").Value = 1.5
$ws.Range("L166").Value = 4
$ws.Range("M166").Value = 5
$ws.Range("N166").Value = 1.5
$ws.Range("O166").Value = 4
$ws.Range("P166").Value = 5
$ws.Range("Q166").Value = -1
$ws.Range("R166").Value = 1.775
$ws.Range("S166").Value = 2.025
$ws.Range("T166").Value = 3
$ws.Range("U166").Value = 1.925
$ws.Range("V166").Value = 1.875
$ws.Range("W166").Value = 0
$ws.Range("X166").Value = 0
$ws.Range("Y166").Value = 0
$ws.Range("Z166").Value = 0
$ws.Range("AA166").Value = 0
